$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Improve Get Tile Image"
$ws.Range("C14").Value = "Done"

$ws.Range("C14").Select()
